# Apply scraped-schedule update for Línea 141 (2026-01-12, scrape @ 06:23:52)
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 06:23:52"
$ws1.Cells.Item(3,1).Value = "Total filas: 54"

$ws1.Cells.Item(28,1).Value = "06:23:52"
$ws1.Cells.Item(28,2).Value = "06:33"
$ws1.Cells.Item(28,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(28,4).Value = 10
$ws1.Cells.Item(28,5).Value = "LP1912"
$ws1.Cells.Item(29,1).Value = "06:23:52"
$ws1.Cells.Item(29,2).Value = "06:43"
$ws1.Cells.Item(29,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(29,4).Value = 20
$ws1.Cells.Item(29,5).Value = "LP1912"
$ws1.Cells.Item(30,1).Value = "04:51:28"
$ws1.Cells.Item(30,2).Value = "06:44"
$ws1.Cells.Item(30,3).Value = "225_C ROCA-H SUR"
$ws1.Cells.Item(30,4).Value = 113
$ws1.Cells.Item(30,5).Value = "LP1912"
$ws1.Cells.Item(31,1).Value = "04:51:28"
$ws1.Cells.Item(31,2).Value = "06:46"
$ws1.Cells.Item(31,3).Value = "215C_EL PATO"
$ws1.Cells.Item(31,4).Value = 115
$ws1.Cells.Item(31,5).Value = "LP1912"
$ws1.Cells.Item(32,1).Value = "05:20:00"
$ws1.Cells.Item(32,2).Value = "06:47"
$ws1.Cells.Item(32,3).Value = "215C_EL PATO"
$ws1.Cells.Item(32,4).Value = 87
$ws1.Cells.Item(32,5).Value = "LP1912"
$ws1.Cells.Item(33,1).Value = "05:51:32"
$ws1.Cells.Item(33,2).Value = "06:59"
$ws1.Cells.Item(33,3).Value = "14_ABASTO"
$ws1.Cells.Item(33,4).Value = 68
$ws1.Cells.Item(33,5).Value = "LP1912"
$ws1.Cells.Item(34,1).Value = "05:20:00"
$ws1.Cells.Item(34,2).Value = "07:00"
$ws1.Cells.Item(34,3).Value = "10_OLMOS"
$ws1.Cells.Item(34,4).Value = 100
$ws1.Cells.Item(34,5).Value = "LP1912"
$ws1.Cells.Item(35,1).Value = "05:20:00"
$ws1.Cells.Item(35,2).Value = "07:00"
$ws1.Cells.Item(35,3).Value = "14_ABASTO"
$ws1.Cells.Item(35,4).Value = 100
$ws1.Cells.Item(35,5).Value = "LP1912"
$ws1.Cells.Item(36,1).Value = "06:23:52"
$ws1.Cells.Item(36,2).Value = "07:01"
$ws1.Cells.Item(36,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(36,4).Value = 38
$ws1.Cells.Item(36,5).Value = "LP1912"
$ws1.Cells.Item(37,1).Value = "05:51:32"
$ws1.Cells.Item(37,2).Value = "07:04"
$ws1.Cells.Item(37,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(37,4).Value = 73
$ws1.Cells.Item(37,5).Value = "LP1912"
$ws1.Cells.Item(38,1).Value = "05:20:00"
$ws1.Cells.Item(38,2).Value = "07:05"
$ws1.Cells.Item(38,3).Value = "15_ABASTO"
$ws1.Cells.Item(38,4).Value = 105
$ws1.Cells.Item(38,5).Value = "LP1912"
$ws1.Cells.Item(39,1).Value = "06:23:52"
$ws1.Cells.Item(39,2).Value = "07:06"
$ws1.Cells.Item(39,3).Value = "225_GOMEZ"
$ws1.Cells.Item(39,4).Value = 43
$ws1.Cells.Item(39,5).Value = "LP1912"
$ws1.Cells.Item(40,1).Value = "05:20:00"
$ws1.Cells.Item(40,2).Value = "07:07"
$ws1.Cells.Item(40,3).Value = "225_GOMEZ"
$ws1.Cells.Item(40,4).Value = 107
$ws1.Cells.Item(40,5).Value = "LP1912"
$ws1.Cells.Item(41,1).Value = "05:51:32"
$ws1.Cells.Item(41,2).Value = "07:11"
$ws1.Cells.Item(41,3).Value = "215A_EL PATO"
$ws1.Cells.Item(41,4).Value = 80
$ws1.Cells.Item(41,5).Value = "LP1912"
$ws1.Cells.Item(42,1).Value = "05:20:00"
$ws1.Cells.Item(42,2).Value = "07:12"
$ws1.Cells.Item(42,3).Value = "215A_EL PATO"
$ws1.Cells.Item(42,4).Value = 112
$ws1.Cells.Item(42,5).Value = "LP1912"
$ws1.Cells.Item(43,1).Value = "05:51:32"
$ws1.Cells.Item(43,2).Value = "07:15"
$ws1.Cells.Item(43,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(43,4).Value = 84
$ws1.Cells.Item(43,5).Value = "LP1912"
$ws1.Cells.Item(44,1).Value = "05:20:00"
$ws1.Cells.Item(44,2).Value = "07:16"
$ws1.Cells.Item(44,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(44,4).Value = 116
$ws1.Cells.Item(44,5).Value = "LP1912"
$ws1.Cells.Item(45,1).Value = "05:51:32"
$ws1.Cells.Item(45,2).Value = "07:21"
$ws1.Cells.Item(45,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(45,4).Value = 90
$ws1.Cells.Item(45,5).Value = "LP1912"
$ws1.Cells.Item(46,1).Value = "06:23:52"
$ws1.Cells.Item(46,2).Value = "07:22"
$ws1.Cells.Item(46,3).Value = "10_OLMOS"
$ws1.Cells.Item(46,4).Value = 59
$ws1.Cells.Item(46,5).Value = "LP1912"
$ws1.Cells.Item(47,1).Value = "05:51:32"
$ws1.Cells.Item(47,2).Value = "07:28"
$ws1.Cells.Item(47,3).Value = "10_OLMOS"
$ws1.Cells.Item(47,4).Value = 97
$ws1.Cells.Item(47,5).Value = "LP1912"
$ws1.Cells.Item(48,1).Value = "05:51:32"
$ws1.Cells.Item(48,2).Value = "07:31"
$ws1.Cells.Item(48,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(48,4).Value = 100
$ws1.Cells.Item(48,5).Value = "LP1912"
$ws1.Cells.Item(49,1).Value = "05:51:32"
$ws1.Cells.Item(49,2).Value = "07:31"
$ws1.Cells.Item(49,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(49,4).Value = 100
$ws1.Cells.Item(49,5).Value = "LP1912"
$ws1.Cells.Item(50,1).Value = "05:51:32"
$ws1.Cells.Item(50,2).Value = "07:32"
$ws1.Cells.Item(50,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(50,4).Value = 101
$ws1.Cells.Item(50,5).Value = "LP1912"
$ws1.Cells.Item(51,1).Value = "05:51:32"
$ws1.Cells.Item(51,2).Value = "07:36"
$ws1.Cells.Item(51,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(51,4).Value = 105
$ws1.Cells.Item(51,5).Value = "LP1912"
$ws1.Cells.Item(52,1).Value = "06:23:52"
$ws1.Cells.Item(52,2).Value = "07:38"
$ws1.Cells.Item(52,3).Value = "10_OLMOS"
$ws1.Cells.Item(52,4).Value = 75
$ws1.Cells.Item(52,5).Value = "LP1912"
$ws1.Cells.Item(53,1).Value = "05:51:32"
$ws1.Cells.Item(53,2).Value = "07:39"
$ws1.Cells.Item(53,3).Value = "10_OLMOS"
$ws1.Cells.Item(53,4).Value = 108
$ws1.Cells.Item(53,5).Value = "LP1912"
$ws1.Cells.Item(54,1).Value = "05:51:32"
$ws1.Cells.Item(54,2).Value = "07:47"
$ws1.Cells.Item(54,3).Value = "14_ABASTO"
$ws1.Cells.Item(54,4).Value = 116
$ws1.Cells.Item(54,5).Value = "LP1912"
$ws1.Cells.Item(55,1).Value = "06:23:52"
$ws1.Cells.Item(55,2).Value = "07:51"
$ws1.Cells.Item(55,3).Value = "215D_EL PATO"
$ws1.Cells.Item(55,4).Value = 88
$ws1.Cells.Item(55,5).Value = "LP1912"
$ws1.Cells.Item(56,1).Value = "06:23:52"
$ws1.Cells.Item(56,2).Value = "08:05"
$ws1.Cells.Item(56,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(56,4).Value = 102
$ws1.Cells.Item(56,5).Value = "LP1912"
$ws1.Cells.Item(57,1).Value = "06:23:52"
$ws1.Cells.Item(57,2).Value = "08:12"
$ws1.Cells.Item(57,3).Value = "15_ABASTO"
$ws1.Cells.Item(57,4).Value = 109
$ws1.Cells.Item(57,5).Value = "LP1912"
$ws1.Cells.Item(58,1).Value = "06:23:52"
$ws1.Cells.Item(58,2).Value = "08:20"
$ws1.Cells.Item(58,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(58,4).Value = 117
$ws1.Cells.Item(58,5).Value = "LP1912"
$ws1.Cells.Item(59,1).Value = "06:23:52"
$ws1.Cells.Item(59,2).Value = "08:22"
$ws1.Cells.Item(59,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(59,4).Value = 119
$ws1.Cells.Item(59,5).Value = "LP1912"

# ---- Sheet: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 06:23:52"
$ws2.Cells.Item(3,1).Value = "Total filas: 11"

$ws2.Cells.Item(16,1).Value = "06:23:52"
$ws2.Cells.Item(16,2).Value = "07:51"
$ws2.Cells.Item(16,3).Value = "215D_EL PATO"
$ws2.Cells.Item(16,4).Value = 88
$ws2.Cells.Item(16,5).Value = "LP1912"

# ---- Sheet: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 06:23:52"
$ws3.Cells.Item(3,1).Value = "Total filas: 9"

$ws3.Cells.Item(14,1).Value = "06:23:52"
$ws3.Cells.Item(14,2).Value = "08:06"
$ws3.Cells.Item(14,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(14,4).Value = 103
$ws3.Cells.Item(14,5).Value = "L6203"
